# Applies the "bugs related to field modification and combo list" update:
#   - Adds two new bug-report rows (7 and 8) to Sheet1
#   - Adjusts the sheet's active view (scrolled down a bit, new selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 (SL = 7): Combo list selection issue
$ws.Range("A8").Value = 7
$ws.Range("F8").Value = "Chosen option should be saved properly"
$ws.Range("G8").Value = "No option is showed as chosen"
$ws.Range("E8").Value = "1. Load any form in mobile                     2. Select any option from the option list                                                                   3. Click 'Save and Exit'    "

# Row 9 (SL = 8): Modification not working during synchronization
$ws.Range("A9").Value = 8
$ws.Range("F9").Value = "Updated data should be showed properly during sychronization in MS Access"
$ws.Range("G9").Value = "Updated data are not showing properly. "
$ws.Range("E9").Value = "1. Load any facility profile in mobile                                                 2. Do any modification in text field and option list                                            3. Click 'Send Now'                                   4. Sync MS Access with Cloud         "

$ws.Range("B8").Value = "Combo list selection is not working (ref: SEC_1.xml)"
$ws.Range("B9").Value = "Modification not works during synchronization (ref: facility profile)"

$ws.Range("I8").Value = "High"
$ws.Range("I9").Value = "High"

# Match style/format of the preceding rows for the two new rows
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A9:Q9").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Row heights specific to the new rows (60pt / 75pt, matching their wrapped content)
$ws.Rows.Item(8).RowHeight = 60
$ws.Rows.Item(9).RowHeight = 75

# Update the sheet view: scroll so row 4 is at the top and select D6
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("D6").Select()
